# PieWithRAPTOR.pptx - "Updates, Small bug fixes"
#
# The chart's callout/leader lines and their percentage labels were nudged
# to better align with the pie wedges they annotate. This touches the two
# freeform leader-line shapes (pl10, pl11) and twelve label textboxes
# (tx12-tx23) that live inside the single group shape on slide 1.
#
# PowerPoint's Shape.Left/Top/Width/Height are Single-precision (float32)
# values expressed in points, while the OOXML stores exact integer EMUs
# (1 pt = 12700 EMU). Naively doing "$emu / 12700" and assigning it back
# loses a little precision once PowerPoint rounds it through a float32 and
# truncates back to EMU, occasionally landing 1 EMU short of the target.
# EmuToPoints nudges the point value up by tiny increments until the
# round-trip lands exactly back on the requested EMU value, so the saved
# file matches the target offsets/extents exactly.

function EmuToPoints {
    param(
        [double]$Emu
    )

    $basePoints = $Emu / 12700.0

    for ($i = 0; $i -lt 400; $i++) {
        $candidate = $basePoints + ($i * 0.0000005)
        $roundTripped = [math]::Floor([double][single]$candidate * 12700.0)
        if ($roundTripped -eq $Emu) {
            return $candidate
        }
    }

    return $basePoints
}

function SetShapePosition {
    param(
        $Shape,
        [double]$X,
        [double]$Y,
        [double]$Cx,
        [double]$Cy
    )

    $Shape.Left = EmuToPoints $X
    $Shape.Top = EmuToPoints $Y
    if ($Cx -ne $null) {
        $Shape.Width = EmuToPoints $Cx
    }
    if ($Cy -ne $null) {
        $Shape.Height = EmuToPoints $Cy
    }
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$grp = $s.Shapes.Item(2)

# Leader line from the "Streptophyta parasite" wedge - shrinks slightly and
# its top-left end moves down/right with the wedge label.
SetShapePosition $grp.GroupItems.Item(8) 5415800 2802222 140323 175655

# Leader line pointing at the neighbouring wedge - grows a little and shifts up.
SetShapePosition $grp.GroupItems.Item(9) 5919287 2920016 248918 133387

# Text labels (percentages / category names) - position only, same size.
SetShapePosition $grp.GroupItems.Item(10) 3340177 2822748
SetShapePosition $grp.GroupItems.Item(11) 3972875 3174865
SetShapePosition $grp.GroupItems.Item(12) 4533739 2309418
SetShapePosition $grp.GroupItems.Item(13) 4889125 2615802
SetShapePosition $grp.GroupItems.Item(14) 6191066 2657951
SetShapePosition $grp.GroupItems.Item(15) 6781418 2964335
SetShapePosition $grp.GroupItems.Item(16) 5955714 3177188
SetShapePosition $grp.GroupItems.Item(17) 6094185 3524118
SetShapePosition $grp.GroupItems.Item(18) 6417309 3735751
SetShapePosition $grp.GroupItems.Item(19) 6591985 4082681
SetShapePosition $grp.GroupItems.Item(20) 3805190 5175123
SetShapePosition $grp.GroupItems.Item(21) 3895492 5522053
